# Update sheet dimensions / values per natmiOut Gas6-Tyro3 LR-pair recompute (Dr Hou advice).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.31916866666667
$ws.Range("H2").Value = 57.957506
$ws.Range("I2").Value = 0.09973928790435696
$ws.Range("J2").Value = 0.1012662650824037
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.2063693333333333
$ws.Range("N2").Value = 0.619108
$ws.Range("O2").Value = 0.05249514260861875
$ws.Range("P2").Value = 0.05463357984752036
$ws.Range("Q2").Value = 3.986883958294222
$ws.Range("R2").Value = 35.881955624648
$ws.Range("S2").Value = 0.005235828142221302
$ws.Range("T2").Value = 0.005532538579239663

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.31916866666667
$ws.Range("H3").Value = 57.957506
$ws.Range("I3").Value = 0.09973928790435696
$ws.Range("J3").Value = 0.1012662650824037
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.258868
$ws.Range("N3").Value = 9.776603999999999
$ws.Range("O3").Value = 0.8289736543672389
$ws.Range("P3").Value = 0.8627426479250582
$ws.Range("Q3").Value = 62.95862055440266
$ws.Range("R3").Value = 566.627584989624
$ws.Range("S3").Value = 0.08268124197806094
$ws.Range("T3").Value = 0.0873667256826738

# Row 4
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.31916866666667
$ws.Range("H4").Value = 57.957506
$ws.Range("I4").Value = 0.09973928790435696
$ws.Range("J4").Value = 0.1012662650824037
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.004351333333333333
$ws.Range("N4").Value = 0.013054
$ws.Range("O4").Value = 0.001106869224130377
$ws.Range("P4").Value = 0.001151958545729551
$ws.Range("Q4").Value = 0.08406414259155555
$ws.Range("R4").Value = 0.756577283324
$ws.Range("S4").Value = 0.0001103983482180118
$ws.Range("T4").Value = 0.0001166545394557889

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.31916866666667
$ws.Range("H5").Value = 57.957506
$ws.Range("I5").Value = 0.09973928790435696
$ws.Range("J5").Value = 0.1012662650824037
$ws.Range("K5").Value = 2
$ws.Range("M5").Value = 0.4616195
$ws.Range("N5").Value = 0.923239
$ws.Range("O5").Value = 0.117424333800012
$ws.Range("P5").Value = 0.0814718136816918
$ws.Range("Q5").Value = 8.918104980322333
$ws.Range("R5").Value = 53.50862988193401
$ws.Range("S5").Value = 0.01171181943585672
$ws.Range("T5").Value = 0.008250346281034404

# Row 6
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 127.3682276666667
$ws.Range("H6").Value = 382.104683
$ws.Range("I6").Value = 0.6575653719009243
$ws.Range("J6").Value = 0.6676324912584373
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.2063693333333333
$ws.Range("N6").Value = 0.619108
$ws.Range("O6").Value = 0.05249514260861875
$ws.Range("P6").Value = 0.05463357984752036
$ws.Range("Q6").Value = 26.28489623141822
$ws.Range("R6").Value = 236.564066082764
$ws.Range("S6").Value = 0.03451898797242845
$ws.Range("T6").Value = 0.03647515301996677

# Row 7
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 127.3682276666667
$ws.Range("H7").Value = 382.104683
$ws.Range("I7").Value = 0.6575653719009243
$ws.Range("J7").Value = 0.6676324912584373
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.258868
$ws.Range("N7").Value = 9.776603999999999
$ws.Range("O7").Value = 0.8289736543672389
$ws.Range("P7").Value = 0.8627426479250582
$ws.Range("Q7").Value = 415.0762413596146
$ws.Range("R7").Value = 3735.686172236532
$ws.Range("S7").Value = 0.5451043693300618
$ws.Range("T7").Value = 0.5759950233491075

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 127.3682276666667
$ws.Range("H8").Value = 382.104683
$ws.Range("I8").Value = 0.6575653719009243
$ws.Range("J8").Value = 0.6676324912584373
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.004351333333333333
$ws.Range("N8").Value = 0.013054
$ws.Range("O8").Value = 0.001106869224130377
$ws.Range("P8").Value = 0.001151958545729551
$ws.Range("Q8").Value = 0.5542216146535555
$ws.Range("R8").Value = 4.987994531882
$ws.Range("S8").Value = 0.0007278388730109788
$ws.Range("T8").Value = 0.0007690849537118664

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 127.3682276666667
$ws.Range("H9").Value = 382.104683
$ws.Range("I9").Value = 0.6575653719009243
$ws.Range("J9").Value = 0.6676324912584373
$ws.Range("K9").Value = 2
$ws.Range("M9").Value = 0.4616195
$ws.Range("N9").Value = 0.923239
$ws.Range("O9").Value = 0.117424333800012
$ws.Range("P9").Value = 0.0814718136816918
$ws.Range("Q9").Value = 58.79565757137284
$ws.Range("R9").Value = 352.773945428237
$ws.Range("S9").Value = 0.0772141757254232
$ws.Range("T9").Value = 0.05439322993565113

# Row 10
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 18.657769
$ws.Range("H10").Value = 55.97330699999999
$ws.Range("I10").Value = 0.09632467245626405
$ws.Range("J10").Value = 0.09779937294404559
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.2063693333333333
$ws.Range("N10").Value = 0.619108
$ws.Range("O10").Value = 0.05249514260861875
$ws.Range("P10").Value = 0.05463357984752036
$ws.Range("Q10").Value = 3.850391350017333
$ws.Range("R10").Value = 34.65352215015599
$ws.Range("S10").Value = 0.005056577417320072
$ws.Range("T10").Value = 0.005343129850775937

# Row 11
$ws.Range("A11").Value = "M1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 18.657769
$ws.Range("H11").Value = 55.97330699999999
$ws.Range("I11").Value = 0.09632467245626405
$ws.Range("J11").Value = 0.09779937294404559
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.258868
$ws.Range("N11").Value = 9.776603999999999
$ws.Range("O11").Value = 0.8289736543672389
$ws.Range("P11").Value = 0.8627426479250582
$ws.Range("Q11").Value = 60.80320634549199
$ws.Range("R11").Value = 547.2288571094278
$ws.Range("S11").Value = 0.07985061573179654
$ws.Range("T11").Value = 0.08437568997915618

# Row 12
$ws.Range("A12").Value = "M1"
$ws.Range("D12").Value = "M1"
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 18.657769
$ws.Range("H12").Value = 55.97330699999999
$ws.Range("I12").Value = 0.09632467245626405
$ws.Range("J12").Value = 0.09779937294404559
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.004351333333333333
$ws.Range("N12").Value = 0.013054
$ws.Range("O12").Value = 0.001106869224130377
$ws.Range("P12").Value = 0.001151958545729551
$ws.Range("Q12").Value = 0.08118617217533332
$ws.Range("R12").Value = 0.7306755495779998
$ws.Range("S12").Value = 0.0001066188154662777
$ws.Range("T12").Value = 0.0001126608234298847

# Row 13
$ws.Range("A13").Value = "M1"
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 18.657769
$ws.Range("H13").Value = 55.97330699999999
$ws.Range("I13").Value = 0.09632467245626405
$ws.Range("J13").Value = 0.09779937294404559
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 0.4616195
$ws.Range("N13").Value = 0.923239
$ws.Range("O13").Value = 0.117424333800012
$ws.Range("P13").Value = 0.0814718136816918
$ws.Range("Q13").Value = 8.612789996895499
$ws.Range("R13").Value = 51.676739981373
$ws.Range("S13").Value = 0.01131086049168118
$ws.Range("T13").Value = 0.007967892290683572

# Row 14
$ws.Range("A14").Value = "M2"
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 19.589352
$ws.Range("H14").Value = 58.768056
$ws.Range("I14").Value = 0.1011341664177781
$ws.Range("J14").Value = 0.1026824987478506
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.2063693333333333
$ws.Range("N14").Value = 0.619108
$ws.Range("O14").Value = 0.05249514260861875
$ws.Range("P14").Value = 0.05463357984752036
$ws.Range("Q14").Value = 4.042641512672
$ws.Range("R14").Value = 36.383773614048
$ws.Range("S14").Value = 0.005309052488705044
$ws.Range("T14").Value = 0.005609912494283605

# Row 15
$ws.Range("A15").Value = "M2"
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 19.589352
$ws.Range("H15").Value = 58.768056
$ws.Range("I15").Value = 0.1011341664177781
$ws.Range("J15").Value = 0.1026824987478506
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.258868
$ws.Range("N15").Value = 9.776603999999999
$ws.Range("O15").Value = 0.8289736543672389
$ws.Range("P15").Value = 0.8627426479250582
$ws.Range("Q15").Value = 63.839112373536
$ws.Range("R15").Value = 574.5520113618239
$ws.Range("S15").Value = 0.08383755951673001
$ws.Range("T15").Value = 0.0885885708652821

# Row 16
$ws.Range("A16").Value = "M2"
$ws.Range("D16").Value = "M1"
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 19.589352
$ws.Range("H16").Value = 58.768056
$ws.Range("I16").Value = 0.1011341664177781
$ws.Range("J16").Value = 0.1026824987478506
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.004351333333333333
$ws.Range("N16").Value = 0.013054
$ws.Range("O16").Value = 0.001106869224130377
$ws.Range("P16").Value = 0.001151958545729551
$ws.Range("Q16").Value = 0.085239800336
$ws.Range("R16").Value = 0.767158203024
$ws.Range("S16").Value = 0.0001119422963159185
$ws.Range("T16").Value = 0.0001182859819294504

# Row 17
$ws.Range("A17").Value = "M2"
$ws.Range("B17").Value = "Gas6"
$ws.Range("C17").Value = "Tyro3"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 19.589352
$ws.Range("H17").Value = 58.768056
$ws.Range("I17").Value = 0.1011341664177781
$ws.Range("J17").Value = 0.1026824987478506
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.4616195
$ws.Range("N17").Value = 0.923239
$ws.Range("O17").Value = 0.117424333800012
$ws.Range("P17").Value = 0.0814718136816918
$ws.Range("Q17").Value = 9.042826875564002
$ws.Range("R17").Value = 54.25696125338401
$ws.Range("S17").Value = 0.01187561211602715
$ws.Range("T17").Value = 0.008365729406355436

# Row 18
$ws.Range("A18").Value = "sCs"
$ws.Range("B18").Value = "Gas6"
$ws.Range("C18").Value = "Tyro3"
$ws.Range("D18").Value = "ECs"
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 8.76216
$ws.Range("H18").Value = 17.52432
$ws.Range("I18").Value = 0.04523650132067659
$ws.Range("J18").Value = 0.03061937196726285
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 0.2063693333333333
$ws.Range("N18").Value = 0.619108
$ws.Range("O18").Value = 0.05249514260861875
$ws.Range("P18").Value = 0.05463357984752036
$ws.Range("Q18").Value = 1.80824111776
$ws.Range("R18").Value = 10.84944670656
$ws.Range("S18").Value = 0.002374696587943888
$ws.Range("T18").Value = 0.001672845903254381

# Row 19
$ws.Range("A19").Value = "sCs"
$ws.Range("B19").Value = "Gas6"
$ws.Range("C19").Value = "Tyro3"
$ws.Range("D19").Value = "FAPs"
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 8.76216
$ws.Range("H19").Value = 17.52432
$ws.Range("I19").Value = 0.04523650132067659
$ws.Range("J19").Value = 0.03061937196726285
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 3.258868
$ws.Range("N19").Value = 9.776603999999999
$ws.Range("O19").Value = 0.8289736543672389
$ws.Range("P19").Value = 0.8627426479250582
$ws.Range("Q19").Value = 28.55472283488
$ws.Range("R19").Value = 171.32833700928
$ws.Range("S19").Value = 0.0374998678105897
$ws.Range("T19").Value = 0.02641663804883864

# Row 20
$ws.Range("A20").Value = "sCs"
$ws.Range("B20").Value = "Gas6"
$ws.Range("C20").Value = "Tyro3"
$ws.Range("D20").Value = "M1"
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 8.76216
$ws.Range("H20").Value = 17.52432
$ws.Range("I20").Value = 0.04523650132067659
$ws.Range("J20").Value = 0.03061937196726285
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.004351333333333333
$ws.Range("N20").Value = 0.013054
$ws.Range("O20").Value = 0.001106869224130377
$ws.Range("P20").Value = 0.001151958545729551
$ws.Range("Q20").Value = 0.03812707888
$ws.Range("R20").Value = 0.22876247328
$ws.Range("S20").Value = 0.00005007089111919006
$ws.Range("T20").Value = 0.00003527224720256028

# Row 21
$ws.Range("A21").Value = "sCs"
$ws.Range("B21").Value = "Gas6"
$ws.Range("C21").Value = "Tyro3"
$ws.Range("D21").Value = "sCs"
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 8.76216
$ws.Range("H21").Value = 17.52432
$ws.Range("I21").Value = 0.04523650132067659
$ws.Range("J21").Value = 0.03061937196726285
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.4616195
$ws.Range("N21").Value = 0.923239
$ws.Range("O21").Value = 0.117424333800012
$ws.Range("P21").Value = 0.0814718136816918
$ws.Range("Q21").Value = 4.04478391812
$ws.Range("R21").Value = 16.17913567248
$ws.Range("S21").Value = 0.005311866031023814
$ws.Range("T21").Value = 0.002494615767967256

